$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15 (hunk 0)
$ws.Range("H15").Value = 2518.1382
$ws.Range("I15").Value = 2518.1382
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 7554.4146
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -7385.4146
# Row 138 (hunk 1)
$ws.Range("H138").Value = 581466.9
$ws.Range("I138").Value = 2155.5715
$ws.Range("J138").Value = 641101.9
$ws.Range("K138").Value = 6466.7145
$ws.Range("L138").Value = 1923305.7
$ws.Range("M138").Value = -1326.7145
$ws.Range("N138").Value = -1933585.7

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (hunk 2)
$ws.Range("H32").Value = 2526.6
$ws.Range("I32").Value = 2019.4884
$ws.Range("J32").Value = 5641.7144
$ws.Range("K32").Value = 2019.4884
$ws.Range("L32").Value = 5641.7144
$ws.Range("M32").Value = -1732.4884
$ws.Range("N32").Value = -6215.7144
# Row 52 (hunk 3)
$ws.Range("H52").Value = 32450
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 32450
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 32450
$ws.Range("N52").Value = -33086
# Row 61 (hunk 4)
$ws.Range("H61").Value = 200001310
$ws.Range("I61").Value = 250000900
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 250000900
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -250000688
$ws.Range("N61").Value = -3424
# Row 74 (hunk 5)
$ws.Range("H74").Value = 2018
$ws.Range("I74").Value = 1542.5714
$ws.Range("J74").Value = 2850
$ws.Range("K74").Value = 1542.5714
$ws.Range("L74").Value = 2850
$ws.Range("M74").Value = -668.5714
$ws.Range("N74").Value = -4598
# Row 77 (hunk 6)
$ws.Range("H77").Value = 2018
$ws.Range("I77").Value = 1542.5714
$ws.Range("J77").Value = 2850
$ws.Range("K77").Value = 7712.857
$ws.Range("L77").Value = 14250
$ws.Range("M77").Value = -3344.857
$ws.Range("N77").Value = -22986
# Row 122 (hunk 7)
$ws.Range("H122").Value = 2515.8
$ws.Range("I122").Value = 2100
$ws.Range("J122").Value = 3139.5
$ws.Range("K122").Value = 6300
$ws.Range("L122").Value = 9418.5
$ws.Range("M122").Value = -3850
$ws.Range("N122").Value = -14318.5
# Row 132 (hunk 8)
$ws.Range("H132").Value = 2067.1555
$ws.Range("I132").Value = 1629.3667
$ws.Range("J132").Value = 2942.7334
$ws.Range("K132").Value = 4888.1001
$ws.Range("L132").Value = 8828.200199999999
$ws.Range("M132").Value = -2358.1001
$ws.Range("N132").Value = -13888.2002
# Row 136 (hunk 9)
$ws.Range("H136").Value = 200001310
$ws.Range("I136").Value = 250000900
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 750002700
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -750000150
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
# Row 13 (hunk 10)
$ws.Range("H13").Value = 32000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 32000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 32000
$ws.Range("N13").Value = -32336
# Row 86 (hunk 11)
$ws.Range("H86").Value = 2633.125
$ws.Range("I86").Value = 2458.7646
$ws.Range("J86").Value = 3056.5715
$ws.Range("K86").Value = 2458.7646
$ws.Range("L86").Value = 3056.5715
$ws.Range("M86").Value = -1335.7646
$ws.Range("N86").Value = -5302.5715
# Row 89 (hunk 12)
$ws.Range("H89").Value = 2633.125
$ws.Range("I89").Value = 2458.7646
$ws.Range("J89").Value = 3056.5715
$ws.Range("K89").Value = 12293.823
$ws.Range("L89").Value = 15282.8575
$ws.Range("M89").Value = -6677.823
$ws.Range("N89").Value = -26514.8575
# Row 134 (hunk 13)
$ws.Range("H134").Value = 1507.8572
$ws.Range("I134").Value = 1211.5
$ws.Range("J134").Value = 2248.75
$ws.Range("K134").Value = 3634.5
$ws.Range("L134").Value = 6746.25
$ws.Range("M134").Value = -1099.5
$ws.Range("N134").Value = -11816.25

$ws = $wb.Worksheets.Item("CRP")
# Row 94 (hunk 14)
$ws.Range("H94").Value = 4269.6665
$ws.Range("I94").Value = 1399
$ws.Range("J94").Value = 5705
$ws.Range("K94").Value = 1399
$ws.Range("L94").Value = 5705
$ws.Range("M94").Value = -948
$ws.Range("N94").Value = -6607
# Row 107 (hunk 15)
$ws.Range("H107").Value = 1025.0769
$ws.Range("I107").Value = 593.25
$ws.Range("J107").Value = 1716
$ws.Range("K107").Value = 593.25
$ws.Range("L107").Value = 1716
$ws.Range("M107").Value = 1326.75
$ws.Range("N107").Value = -5556
# Row 132 (hunk 16)
$ws.Range("H132").Value = 1232.0968
$ws.Range("I132").Value = 893.381
$ws.Range("J132").Value = 1943.4
$ws.Range("K132").Value = 2680.143
$ws.Range("L132").Value = 5830.200000000001
$ws.Range("M132").Value = -150.143
$ws.Range("N132").Value = -10890.2
# Row 134 (hunk 17)
$ws.Range("H134").Value = 31252666
$ws.Range("I134").Value = 3207.6667
$ws.Range("J134").Value = 71430540
$ws.Range("K134").Value = 9623.000100000001
$ws.Range("L134").Value = 214291620
$ws.Range("M134").Value = -7088.000100000001
$ws.Range("N134").Value = -214296690

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (hunk 18)
$ws.Range("H5").Value = 1039.4375
$ws.Range("I5").Value = 1178.1305
$ws.Range("J5").Value = 685
$ws.Range("K5").Value = 3534.3915
$ws.Range("L5").Value = 2055
$ws.Range("M5").Value = -3422.3915
$ws.Range("N5").Value = -2279
# Row 74 (hunk 19)
$ws.Range("H74").Value = 4493.647
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 4493.647
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 13480.941
$ws.Range("N74").Value = -15602.941
$ws.Range("M74").ClearContents()
# Row 77 (hunk 20)
$ws.Range("H77").Value = 4493.647
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 4493.647
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 40442.823
$ws.Range("N77").Value = -51050.823
$ws.Range("M77").ClearContents()
# Row 93 (hunk 21)
$ws.Range("H93").Value = 7514.5
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 7514.5
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 22543.5
$ws.Range("N93").Value = -26287.5
# Row 131 (hunk 22)
$ws.Range("H131").Value = 12987881
$ws.Range("I131").Value = 55556040
$ws.Range("J131").Value = 984.678
$ws.Range("K131").Value = 166668120
$ws.Range("L131").Value = 2954.034
$ws.Range("M131").Value = -166663080
$ws.Range("N131").Value = -13034.034
# Row 135 (hunk 23)
$ws.Range("H135").Value = 1039.4375
$ws.Range("I135").Value = 1178.1305
$ws.Range("J135").Value = 685
$ws.Range("K135").Value = 10603.1745
$ws.Range("L135").Value = 6165
$ws.Range("M135").Value = -8068.174499999999
$ws.Range("N135").Value = -11235

$ws = $wb.Worksheets.Item("GSM")
# Row 102 (hunk 24)
$ws.Range("H102").Value = 1143.7894
$ws.Range("I102").Value = 1218.4
$ws.Range("J102").Value = 1060.8889
$ws.Range("K102").Value = 1218.4
$ws.Range("L102").Value = 1060.8889
$ws.Range("M102").Value = 403.5999999999999
$ws.Range("N102").Value = -4304.8889
# Row 116 (hunk 25)
$ws.Range("H116").Value = 35000
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 35000
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 35000
$ws.Range("N116").Value = -44178
# Row 122 (hunk 26)
$ws.Range("H122").Value = 3475.1875
$ws.Range("I122").Value = 3635.25
$ws.Range("J122").Value = 2995
$ws.Range("K122").Value = 10905.75
$ws.Range("L122").Value = 8985
$ws.Range("M122").Value = -8455.75
$ws.Range("N122").Value = -13885
# Row 132 (hunk 27)
$ws.Range("H132").Value = 3643.8262
$ws.Range("I132").Value = 3848.5833
$ws.Range("J132").Value = 3420.4546
$ws.Range("K132").Value = 11545.7499
$ws.Range("L132").Value = 10261.3638
$ws.Range("M132").Value = -9015.749899999999
$ws.Range("N132").Value = -15321.3638

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (hunk 28)
$ws.Range("H7").Value = 2054.375
$ws.Range("I7").Value = 2006.0834
$ws.Range("J7").Value = 2199.25
$ws.Range("K7").Value = 2006.0834
$ws.Range("L7").Value = 2199.25
$ws.Range("M7").Value = -1894.0834
$ws.Range("N7").Value = -2423.25
# Row 100 (hunk 29)
$ws.Range("H100").Value = 1215.2727
$ws.Range("I100").Value = 1081.1428
$ws.Range("J100").Value = 1450
$ws.Range("K100").Value = 1081.1428
$ws.Range("L100").Value = 1450
$ws.Range("M100").Value = -540.1428000000001
$ws.Range("N100").Value = -2532
# Row 123 (hunk 30)
$ws.Range("H123").Value = 40959
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 40959
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 40959
$ws.Range("N123").Value = -50759
# Row 126 (hunk 31)
$ws.Range("H126").Value = 2054.375
$ws.Range("I126").Value = 2006.0834
$ws.Range("J126").Value = 2199.25
$ws.Range("K126").Value = 6018.2502
$ws.Range("L126").Value = 6597.75
$ws.Range("M126").Value = -3548.2502
$ws.Range("N126").Value = -11537.75
# Row 136 (hunk 32)
$ws.Range("H136").Value = 1434.25
$ws.Range("I136").Value = 1080.2307
$ws.Range("J136").Value = 2968.3333
$ws.Range("K136").Value = 3240.6921
$ws.Range("L136").Value = 8904.999899999999
$ws.Range("M136").Value = -690.6921000000002
$ws.Range("N136").Value = -14004.9999

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (hunk 33)
$ws.Range("H81").Value = 344.66666
$ws.Range("I81").Value = 500
$ws.Range("J81").Value = 267
$ws.Range("K81").Value = 1000
$ws.Range("L81").Value = 534
$ws.Range("M81").Value = 61
$ws.Range("N81").Value = -2656
# Row 84 (hunk 34)
$ws.Range("H84").Value = 344.66666
$ws.Range("I84").Value = 500
$ws.Range("J84").Value = 267
$ws.Range("K84").Value = 5000
$ws.Range("L84").Value = 2670
$ws.Range("M84").Value = 304
$ws.Range("N84").Value = -13278
# Row 122 (hunk 35)
$ws.Range("H122").Value = 20834898
$ws.Range("I122").Value = 22728798
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 68186394
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -68183944
$ws.Range("N122").Value = -10900
# Row 126 (hunk 36)
$ws.Range("H126").Value = 71433550
$ws.Range("I126").Value = 90912070
$ws.Range("J126").Value = 12333.333
$ws.Range("K126").Value = 272736210
$ws.Range("L126").Value = 36999.999
$ws.Range("M126").Value = -272733740
$ws.Range("N126").Value = -41939.999
# Row 141 (hunk 37)
$ws.Range("H141").Value = 39815.555
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 39815.555
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 39815.555
$ws.Range("N141").Value = -50175.555
